$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new date cells (column A) use the same date style as existing data (s="2" -> custom date format)
$dateFormat = $ws.Range("A1516").NumberFormat
$ws.Range("A1517:A1540").NumberFormat = $dateFormat

# Row 1517
$ws.Range("A1517").Value = 44459
$ws.Range("B1517").Value = "Reko272"
$ws.Range("C1517").Value = 3011
$ws.Range("D1517").Value = "Reko Swish +46733987105"
$ws.Range("F1517").Value = 141.07

# Row 1518
$ws.Range("A1518").Value = 44459
$ws.Range("B1518").Value = "Reko272"
$ws.Range("C1518").Value = 2611
$ws.Range("D1518").Value = "Reko Swish +46733987105"
$ws.Range("F1518").Value = 16.93

# Row 1519
$ws.Range("A1519").Value = 44459
$ws.Range("B1519").Value = "Reko272"
$ws.Range("C1519").Value = 1930
$ws.Range("D1519").Value = "Reko Swish +46733987105"
$ws.Range("E1519").Value = 158

# Row 1520
$ws.Range("A1520").Value = 44461
$ws.Range("B1520").Value = "Reko273"
$ws.Range("C1520").Value = 3011
$ws.Range("D1520").Value = "Reko Swish +46723698764"
$ws.Range("F1520").Value = 230.36

# Row 1521
$ws.Range("A1521").Value = 44461
$ws.Range("B1521").Value = "Reko273"
$ws.Range("C1521").Value = 2611
$ws.Range("D1521").Value = "Reko Swish +46723698764"
$ws.Range("F1521").Value = 27.64

# Row 1522
$ws.Range("A1522").Value = 44461
$ws.Range("B1522").Value = "Reko273"
$ws.Range("C1522").Value = 1930
$ws.Range("D1522").Value = "Reko Swish +46723698764"
$ws.Range("E1522").Value = 258

# Row 1523
$ws.Range("A1523").Value = 44461
$ws.Range("B1523").Value = "Reko274"
$ws.Range("C1523").Value = 3011
$ws.Range("D1523").Value = "Reko Swish +46709906521"
$ws.Range("F1523").Value = 230.36

# Row 1524
$ws.Range("A1524").Value = 44461
$ws.Range("B1524").Value = "Reko274"
$ws.Range("C1524").Value = 2611
$ws.Range("D1524").Value = "Reko Swish +46709906521"
$ws.Range("F1524").Value = 27.64

# Row 1525
$ws.Range("A1525").Value = 44461
$ws.Range("B1525").Value = "Reko274"
$ws.Range("C1525").Value = 1930
$ws.Range("D1525").Value = "Reko Swish +46709906521"
$ws.Range("E1525").Value = 258

# Row 1526
$ws.Range("A1526").Value = 44461
$ws.Range("B1526").Value = "Reko275"
$ws.Range("C1526").Value = 3011
$ws.Range("D1526").Value = "Reko Swish +46703677212"
$ws.Range("F1526").Value = 1062.5

# Row 1527
$ws.Range("A1527").Value = 44461
$ws.Range("B1527").Value = "Reko275"
$ws.Range("C1527").Value = 2611
$ws.Range("D1527").Value = "Reko Swish +46703677212"
$ws.Range("F1527").Value = 127.5

# Row 1528
$ws.Range("A1528").Value = 44461
$ws.Range("B1528").Value = "Reko275"
$ws.Range("C1528").Value = 1930
$ws.Range("D1528").Value = "Reko Swish +46703677212"
$ws.Range("E1528").Value = 1190

# Row 1529
$ws.Range("A1529").Value = 44462
$ws.Range("C1529").Value = 6540
$ws.Range("D1529").Value = "Klarna*kamda.se K0135"
$ws.Range("E1529").Value = 479.2

# Row 1530
$ws.Range("A1530").Value = 44462
$ws.Range("C1530").Value = 2641
$ws.Range("D1530").Value = "Klarna*kamda.se K0135"
$ws.Range("E1530").Value = 119.8

# Row 1531
$ws.Range("A1531").Value = 44462
$ws.Range("C1531").Value = 1930
$ws.Range("D1531").Value = "Klarna*kamda.se K0135"
$ws.Range("F1531").Value = 599

# Row 1532
$ws.Range("A1532").Value = 44464
$ws.Range("C1532").Value = 7010
$ws.Range("D1532").Value = "Lön September"
$ws.Range("E1532").Value = 1317

# Row 1533
$ws.Range("A1533").Value = 44464
$ws.Range("D1533").Value = "Lön September"
$ws.Range("E1533").Value = 0

# Row 1534
$ws.Range("A1534").Value = 44464
$ws.Range("C1534").Value = 1930
$ws.Range("D1534").Value = "Lön September"
$ws.Range("F1534").Value = 1317

# Row 1535
$ws.Range("A1535").Value = 44464
$ws.Range("C1535").Value = 6400
$ws.Range("D1535").Value = "FACEBK FG24667Z62 K0135"
$ws.Range("E1535").Value = 100

# Row 1536
$ws.Range("A1536").Value = 44464
$ws.Range("D1536").Value = "FACEBK FG24667Z62 K0135"
$ws.Range("E1536").Value = 0

# Row 1537
$ws.Range("A1537").Value = 44464
$ws.Range("C1537").Value = 1930
$ws.Range("D1537").Value = "FACEBK FG24667Z62 K0135"
$ws.Range("F1537").Value = 100

# Row 1538
$ws.Range("A1538").Value = 44465
$ws.Range("C1538").Value = 7010
$ws.Range("D1538").Value = "Sinthu lön Septe"
$ws.Range("E1538").Value = 4584

# Row 1539
$ws.Range("A1539").Value = 44465
$ws.Range("D1539").Value = "Sinthu lön Septe"
$ws.Range("E1539").Value = 0

# Row 1540
$ws.Range("A1540").Value = 44465
$ws.Range("C1540").Value = 1930
$ws.Range("D1540").Value = "Sinthu lön Septe"
$ws.Range("F1540").Value = 4584
